$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update H12 (Sunday hours for week commencing 2018-04-29) from 3 to 6.
# The dependent formulas in I12 (row total) and I19 (grand total) recalc automatically.
$ws.Range("H12").Value = 6

# Update the active selection to match the saved view state (M10).
$ws.Activate()
$ws.Range("M10").Select()

$excel.Calculate()

$wb.Save()
